$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.786498546600342
$ws.Range("B1").Value = 4.403322696685791
$ws.Range("C1").Value = 3.117683172225952
$ws.Range("D1").Value = 2.429234266281128
$ws.Range("E1").Value = 2.139770746231079
